$d = $word.ActiveDocument

function ReplaceIn($range, $findText, $replaceText) {
    $ok = $range.Find.Execute(
        $findText, $true, $false, $false, $false, $false,
        $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $findText"
    }
}

# ===========================================================================
# Paragraph 1 - English blog post body
# ===========================================================================
$p1 = $d.Paragraphs(1).Range

# "I don't like the story I told back then," -> "I don't like the old storyboards"
# (this also removes the gramStart/gramEnd proofErr pair that wrapped "then,")
ReplaceIn $p1 `
    "I don’t like the story I told back then, I think my storytelling style was not quite there yet. " `
    "I don’t like the old storyboards, I think my storytelling style was not quite there yet. "

# New sentence about redrawing old comics that were already posted
ReplaceIn $p1 `
    "Since the last failed attempt at an English-Portuguese version of the website, I broke a bunch of stuff and apparently the buttons are not working on iPhones. " `
    "I’ve been redrawing some of the comics that were already posted (looking at you 1 and 5), because I drew some of them in a hurry. "

# New paragraph about the website / React / Javascript Vanilla rewrite, replacing
# the old "cache / buttons / journey" sentences and the old "Soon I'll recode..." run
ReplaceIn $p1 `
    "Also, for some reason the cache doesn’t update on mobile so when people access the website first thing Monday morning, they still see the previous comic. It’s been a crazy journey guys, but I couldn’t ask for a better prototype to test my coding skills. Soon I’ll recode everything using the framework React. Bear with me" `
    "About the website: look at the amazing job I did putting up Portuguese and English options with a button and all! (LIES – I asked the frontend senior that sits by my side to code the entire thing for me). However, he did use Javascript Vanilla cause that’s what I was going for, but I will recode the whole website myself (this time I will do it myself, I promise) in React, so I will probably fuck up some things along the week. Bear with me"

# ===========================================================================
# Paragraph 3 - Portuguese blog post body
# ===========================================================================
$p3 = $d.Paragraphs(3).Range

# "...terceiro personagem no comic." -> "...terceiro personagem no antigo comic."
ReplaceIn $p3 `
    "como um terceiro personagem no comic. Olhando pra trás, eu " `
    "como um terceiro personagem no antigo comic. Olhando pra trás, eu "

# New sentence about redrawing some of the already-posted comics
ReplaceIn $p3 `
    "de contação de história. Sobre o site: minha última tentativa de fazer uma versão em português e inglês gerou milhões de bugs. " `
    "de contação de história. Ah! E eu to redesenhando alguns comics já postados (como o número 1 e 5) porque havia previamente feito na correria. Sobre o site: "

# New paragraph about the website / React / front-end rewrite, replacing the old
# "bugs / cache / buttons / iPhone / javascript journey" sentences
ReplaceIn $p3 `
    "O cacho no mobile não está atualizando, então quem tenta ler o quadrinho segunda cedinho não vê o último quadrinho, e sim o da semana anterior. Os botões pararam de funcionar no iPhone. Essa jornada do meu aprendizado de javascript tá sendo uma montanha-russa, mas é um protótipo perfeito pra testar as coisas que eu estou aprendendo. Em breve vou reescrever o código inteiro em React. Aguenta aí galera" `
    "olhem que maravilhoso trabalho eu fiz colocando um botão pra mudar o idioma (MEN TI RAS – quem fez isso foi o front sênior que trabalha ao meu lado). Ele fez todo código pra mim, mas eu prometo que vou redesenhar tudinho do zero (desta vez, eu mesma farei) usando React. Aguenta aí galera"

Write-Output "done"
